$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the previously-blank "Assignee" cells in column D
$ws.Range("D5").Value  = "gomana"
$ws.Range("D12").Value = "karim"
$ws.Range("D13").Value = "gomana"
$ws.Range("D16").Value = "karim"
$ws.Range("D17").Value = "gomana"
$ws.Range("D20").Value = "karim"
$ws.Range("D21").Value = "gomana"
$ws.Range("D24").Value = "karim"
$ws.Range("D25").Value = "gomana"
$ws.Range("D28").Value = "karim"
$ws.Range("D29").Value = "gomana"
$ws.Range("D32").Value = "islam"
$ws.Range("D33").Value = "islam"
$ws.Range("D36").Value = "islam"
$ws.Range("D37").Value = "islam"
$ws.Range("D39").Value = "nagy"
$ws.Range("D40").Value = "rawi"
$ws.Range("D41").Value = "rawi"

# Update the active selection to match the saved view state
$ws.Range("E16").Select()
